# Natmi following Dr Hou advice
# Expand the Fgf16-Fgfr4 sending/target cluster pairs from 2 rows to the
# full 8-row combination set (FAPs/sCs sending clusters x ECs/FAPs/M2/sCs target clusters),
# refreshing every derived NATMI metric column (E:T) accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colNames = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$rows = @(
    @{ Row = 2; Values = @("FAPs", "Fgf16", "Fgfr4", "ECs", "2", "0.6666666666666666", "0.09755033333333334", "0.292651", "0.1198375636346959", "0.1198375636346959", "3", "1", "0.8372556666666666", "2.511767", "0.04323017762831461", "0.04323017762831462", "0.08167456936855555", "0.7350711243169999", "0.005180599162472359", "0.005180599162472359") },
    @{ Row = 3; Values = @("FAPs", "Fgf16", "Fgfr4", "FAPs", "2", "0.6666666666666666", "0.09755033333333334", "0.292651", "0.1198375636346959", "0.1198375636346959", "2", "0.6666666666666666", "0.110372", "0.331116", "0.005698858013333649", "0.00569885801333365", "0.01076682539066666", "0.09690142851599999", "0.0006829372598179678", "0.0006829372598179679") },
    @{ Row = 4; Values = @("FAPs", "Fgf16", "Fgfr4", "M2", "2", "0.6666666666666666", "0.09755033333333334", "0.292651", "0.1198375636346959", "0.1198375636346959", "1", "0.3333333333333333", "0.01411566666666667", "0.042347", "0.0007288368435552497", "0.0007288368435552499", "0.001376987988555556", "0.012392891897", "0.00008734203161886315", "0.00008734203161886316") },
    @{ Row = 5; Values = @("FAPs", "Fgf16", "Fgfr4", "sCs", "2", "0.6666666666666666", "0.09755033333333334", "0.292651", "0.1198375636346959", "0.1198375636346959", "3", "1", "18.40564566666667", "55.216937", "0.9503421275147964", "0.9503421275147965", "1.795476869998556", "16.159291829987", "0.1138866851807867", "0.1138866851807867") },
    @{ Row = 6; Values = @("sCs", "Fgf16", "Fgfr4", "ECs", "3", "1", "0.716471", "2.149413", "0.8801624363653041", "0.880162436365304", "3", "1", "0.8372556666666666", "2.511767", "0.04323017762831461", "0.04323017762831462", "0.5998694047523333", "5.398824642770999", "0.03804957846584225", "0.03804957846584225") },
    @{ Row = 7; Values = @("sCs", "Fgf16", "Fgfr4", "FAPs", "3", "1", "0.716471", "2.149413", "0.8801624363653041", "0.880162436365304", "2", "0.6666666666666666", "0.110372", "0.331116", "0.005698858013333649", "0.00569885801333365", "0.07907833721199999", "0.711705034908", "0.005015920753515681", "0.005015920753515682") },
    @{ Row = 8; Values = @("sCs", "Fgf16", "Fgfr4", "M2", "3", "1", "0.716471", "2.149413", "0.8801624363653041", "0.880162436365304", "1", "0.3333333333333333", "0.01411566666666667", "0.042347", "0.0007288368435552497", "0.0007288368435552499", "0.01011346581233333", "0.091021192311", "0.0006414948119363866", "0.0006414948119363867") },
    @{ Row = 9; Values = @("sCs", "Fgf16", "Fgfr4", "sCs", "3", "1", "0.716471", "2.149413", "0.8801624363653041", "0.880162436365304", "3", "1", "18.40564566666667", "55.216937", "0.9503421275147964", "0.9503421275147965", "13.18711135644233", "118.684002207981", "0.8364554423340097", "0.8364554423340097") }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $vals = $r.Values
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ref = "{0}{1}" -f $colNames[$i], $rowNum
        $ws.Range($ref).Value = $vals[$i]
    }
}
